# "Classes" sheet: add two new trailing columns, X (CategoriaRvt) and
# Y (ClasseIfc), mirroring the layout/format of the existing "Key" column (V).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")
$ws.Activate()

# Copy the header (row 1) and body (rows 2-19) formatting from column V so the
# two new columns look like the rest of the table (fill/border/font/alignment).
$ws.Range("V1").Copy()
$ws.Range("X1:Y1").PasteSpecial(-4122)

$ws.Range("V2:V19").Copy()
$ws.Range("X2:Y19").PasteSpecial(-4122)

# Header labels.
$ws.Range("X1").Value = "CategoriaRvt"
$ws.Range("Y1").Value = "ClasseIfc"

# Placeholder body values, same "null" text used by the other not-yet-filled
# columns (G:K) in this sheet.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 24).Value = "null"
    $ws.Cells.Item($r, 25).Value = "null"
}

# Leave the new columns selected, matching the author's post-edit selection.
$ws.Range("X2:Y19").Select() | Out-Null
